# Append the 05/29/2020 daily update row to the "Condicion_Pacientes" table,
# matching the author's "Actualizar 05-30-2020" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Excel Table by one row (extends table ref + autofilter, like
# Excel does automatically when you type into the row right below a table).
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

# New day's figures: Fecha, Pruebas Realizadas, Pruebas Positivas,
# Clinicamente Estables, Clinicamente Graves, Cuidados Intensivos.
$ws.Range("A78").Value = 43980
$ws.Range("B78").Value = 567
$ws.Range("C78").Value = 134
$ws.Range("D78").Value = 460
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 34

# Carry the formatting down from the row above (date number format in
# column A, centered alignment in B:F) exactly like Excel's table
# auto-fill-down behavior for a freshly inserted row.
$ws.Range("A77:F77").Copy()
$ws.Range("A78:F78").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection on the last cell of the new row, as in the source file.
[void]$ws.Range("F78").Select()
